$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains its existing Text format so numeric-looking
# price strings (e.g. "1.010") are not reinterpreted as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.319.99'
$ws.Range("E2").Value = '  +1.04%  '

$ws.Range("D3").Value = '1.874.86'
$ws.Range("E3").Value = '  +0.53%  '

$ws.Range("D4").Value = '1.010'
$ws.Range("E4").Value = '  +0.64%  '

$ws.Range("D5").Value = '314.75'
$ws.Range("E5").Value = '  +0.68%  '

$ws.Range("D6").Value = '1.010'
$ws.Range("E6").Value = '  +0.66%  '

$ws.Range("D7").Value = '0.5132'
$ws.Range("E7").Value = '  +0.14%  '

$ws.Range("D8").Value = '0.3932'
$ws.Range("E8").Value = '  +1.75%  '

$ws.Range("D9").Value = '0.08299'
$ws.Range("E9").Value = '  -0.61%  '

$ws.Range("D10").Value = '1.120'
$ws.Range("E10").Value = '  +0.87%  '

$ws.Range("D11").Value = '41.80'
$ws.Range("E11").Value = '  +1.14%  '

$ws.Range("D12").Value = '6.294'
$ws.Range("E12").Value = '  +1.90%  '

$ws.Range("D13").Value = '1.863.13'
$ws.Range("E13").Value = '  -0.38%  '

$ws.Range("D14").Value = '20.36'
$ws.Range("E14").Value = '  -0.38%  '

$ws.Range("D15").Value = '7.248'
$ws.Range("E15").Value = '  -0.43%  '

$ws.Range("D16").Value = '1.010'
$ws.Range("E16").Value = '  +0.63%  '

$ws.Range("D17").Value = '0.00001106'
$ws.Range("E17").Value = '  +0.76%  '

$ws.Range("D18").Value = '91.31'
$ws.Range("E18").Value = '  +0.54%  '

$ws.Range("D19").Value = '0.06723'
$ws.Range("E19").Value = '  +1.31%  '

$ws.Range("D20").Value = '17.74'
$ws.Range("E20").Value = '  +0.76%  '

$ws.Range("D21").Value = '1.009'
$ws.Range("E21").Value = '  +0.60%  '

$ws.Range("D22").Value = '6.002'
$ws.Range("E22").Value = '  +0.30%  '

$ws.Range("D23").Value = '28.341.66'
$ws.Range("E23").Value = '  +1.00%  '

$ws.Range("D24").Value = '11.15'
$ws.Range("E24").Value = '  +1.50%  '

$ws.Range("D25").Value = '2.257'
$ws.Range("E25").Value = '  +0.59%  '

$ws.Range("D26").Value = '2.083.46'
$ws.Range("E26").Value = '  +0.35%  '

$ws.Range("D27").Value = '160.44'
$ws.Range("E27").Value = '  +1.54%  '

$ws.Range("E28").Value = '  +1.42%  '

$ws.Range("D29").Value = '2.443'
$ws.Range("E29").Value = '  -0.66%  '

$ws.Range("D30").Value = '126.80'
$ws.Range("E30").Value = '  +1.57%  '

$ws.Range("D31").Value = '0.1060'
$ws.Range("E31").Value = '  -0.12%  '

$ws.Range("D32").Value = '1.047'
$ws.Range("E32").Value = '  +1.96%  '

$ws.Range("D33").Value = '5.896'
$ws.Range("E33").Value = '  +1.53%  '

$ws.Range("D34").Value = '3.617'
$ws.Range("E34").Value = '  +0.59%  '

$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D35").Value = '0.02441'
$ws.Range("E35").Value = '  +1.03%  '

$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").Value = '9.277'
$ws.Range("E36").Value = '  -1.59%  '

$ws.Range("D37").Value = '0.06548'
$ws.Range("E37").Value = '  +0.44%  '

$ws.Range("D38").Value = '0.2188'
$ws.Range("E38").Value = '  +0.71%  '

$ws.Range("D39").Value = '0.6475'
$ws.Range("E39").Value = '  +0.20%  '

$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '1.246'
$ws.Range("E40").Value = '  +2.90%  '

$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").Value = '1.185'
$ws.Range("E41").Value = '  -1.09%  '

$ws.Range("D42").Value = '4.977'
$ws.Range("E42").Value = '  +0.12%  '

$ws.Range("D43").Value = '11.17'
$ws.Range("E43").Value = '  -0.82%  '

$ws.Range("D44").Value = '0.6070'
$ws.Range("E44").Value = '  +0.07%  '

$ws.Range("D45").Value = '13.12'
$ws.Range("E45").Value = '  +0.52%  '

$ws.Range("D46").Value = '3.693'
$ws.Range("E46").Value = '  +0.64%  '

$ws.Range("D47").Value = '1.277'
$ws.Range("E47").Value = '  -0.81%  '

$ws.Range("D48").Value = '2.011'
$ws.Range("E48").Value = '  +0.70%  '

$ws.Range("D49").Value = '1.218'
$ws.Range("E49").Value = '  -0.04%  '

$ws.Range("D50").Value = '121.87'
$ws.Range("E50").Value = '  +1.16%  '

$ws.Range("D51").Value = '0.06887'
$ws.Range("E51").Value = '  +0.48%  '
